$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31.466290000000001
$ws.Range("B3").Value = 31.843699999999998
$ws.Range("B4").Value = 34.237400000000001
$ws.Range("B6").Value = 40.640720000000002
$ws.Range("B7").Value = 54.841230000000003
$ws.Range("B8").Value = 44.170110000000001
$ws.Range("B9").Value = 33.823740000000001
$ws.Range("B10").Value = 47.657049999999998
$ws.Range("B11").Value = 48.421579999999999
$ws.Range("B12").Value = 62.345860000000002
$ws.Range("B14").Value = 68.945639999999997
$ws.Range("B15").Value = 69.273349999999994

$ws.Range("D16").Select()
